$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Master")

# Extend the Table1 ListObject by one column (K -> L) to make room for the
# new "number" column.
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:L103"))

# Give the new header cell the same formatting as the rest of the header
# row (bold / centered / bordered), then set its text.
$ws.Range("K1").Copy($ws.Range("L1"))
$ws.Range("L1").Value = "number"

# Index every card: row 14 (first card with worker data) .. row 103 get
# sequential numbers 1..90.
for ($r = 14; $r -le 103; $r++) {
    $ws.Cells.Item($r, 12).Value = $r - 13
}

# Reflect the user's selection of the newly filled column.
$ws.Range("L14:L103").Select()
